$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A242').Value = '(Intercept)'
$ws.Range('B242').Value = [double]"2.12022614487124e-86"
$ws.Range('C242').Value = 340.3071951568477
$ws.Range('D242').Value = -0.5796844675924698
$ws.Range('E242').Value = 0.562127419664201
$ws.Range('F242').Value = 0
$ws.Range('G242').Value = [double]"1.148545122303125e+204"
$ws.Range('H242').Value = 'TZP'
$ws.Range('A243').Value = 'Year'
$ws.Range('B243').Value = 1.103265964947725
$ws.Range('C243').Value = 0.1682942745768018
$ws.Range('D243').Value = 0.5839464244642122
$ws.Range('E243').Value = 0.5592563577839413
$ws.Range('F243').Value = 0.7932263952926472
$ws.Range('G243').Value = 1.537202823979058
$ws.Range('H243').Value = 'TZP'
$ws.Range('A244').Value = 'Specimen_typeother'
$ws.Range('B244').Value = 0.4827938278208103
$ws.Range('C244').Value = 0.3118814447319059
$ws.Range('D244').Value = -2.334751189186488
$ws.Range('E244').Value = 0.01955642242869515
$ws.Range('F244').Value = 0.2579391837868429
$ws.Range('G244').Value = 0.8795302289253522
$ws.Range('H244').Value = 'TZP'
$ws.Range('A245').Value = 'Specimen_typeRespiratory'
$ws.Range('B245').Value = 0.4087218536795503
$ws.Range('C245').Value = 0.2857244209323956
$ws.Range('D245').Value = -3.131410383898778
$ws.Range('E245').Value = 0.001739688717862533
$ws.Range('F245').Value = 0.2310474285241748
$ws.Range('G245').Value = 0.7100213466474574
$ws.Range('H245').Value = 'TZP'
$ws.Range('A246').Value = 'Specimen_typeUrine'
$ws.Range('B246').Value = 0.4064758692059096
$ws.Range('C246').Value = 0.2199860458947055
$ws.Range('D246').Value = -4.092217351523558
$ws.Range('E246').Value = [double]"4.272678815771846e-05"
$ws.Range('F246').Value = 0.2639391232547947
$ws.Range('G246').Value = 0.6259422126029462
$ws.Range('H246').Value = 'TZP'
$ws.Range('A247').Value = 'Specimen_typeWound & soft tissues'
$ws.Range('B247').Value = 0.5707325276830005
$ws.Range('C247').Value = 0.2975141226203328
$ws.Range('D247').Value = -1.885068856593717
$ws.Range('E247').Value = 0.05942054206724073
$ws.Range('F247').Value = 0.3157064033360731
$ws.Range('G247').Value = 1.016127416490391
$ws.Range('H247').Value = 'TZP'
$ws.Range('A248').Value = 'HospitalCHBH'
$ws.Range('B248').Value = 0.1490956686064201
$ws.Range('C248').Value = 0.4836803244361255
$ws.Range('D248').Value = -3.934762304236533
$ws.Range('E248').Value = [double]"8.327907350307602e-05"
$ws.Range('F248').Value = 0.05674518583748864
$ws.Range('G248').Value = 0.3795443056860292
$ws.Range('H248').Value = 'TZP'
$ws.Range('A249').Value = 'HospitalCNGMO'
$ws.Range('B249').Value = 0.5490196302122109
$ws.Range('C249').Value = 0.6222527291519727
$ws.Range('D249').Value = -0.9636294928310511
$ws.Range('E249').Value = 0.3352317149440498
$ws.Range('F249').Value = 0.1426669015145889
$ws.Range('G249').Value = 1.725717958582835
$ws.Range('H249').Value = 'TZP'
$ws.Range('A250').Value = 'Ward_ED_ICUED'
$ws.Range('B250').Value = 0.1475795473771434
$ws.Range('C250').Value = 0.4294022124291951
$ws.Range('D250').Value = -4.455934061211244
$ws.Range('E250').Value = [double]"8.352870499896662e-06"
$ws.Range('F250').Value = 0.06217671209713808
$ws.Range('G250').Value = 0.3362371932977084
$ws.Range('H250').Value = 'TZP'
$ws.Range('A251').Value = 'Ward_ED_ICUOther'
$ws.Range('B251').Value = 0.208163910641488
$ws.Range('C251').Value = 0.3469777575472706
$ws.Range('D251').Value = -4.523141450849089
$ws.Range('E251').Value = [double]"6.092849177486386e-06"
$ws.Range('F251').Value = 0.103335142314097
$ws.Range('G251').Value = 0.4049700865509483
$ws.Range('H251').Value = 'TZP'
$ws.Range('A252').Value = 'GenderF'
$ws.Range('B252').Value = 1.009540986485855
$ws.Range('C252').Value = 0.159033548788717
$ws.Range('D252').Value = 0.05970915443534012
$ws.Range('E252').Value = 0.9523872805609095
$ws.Range('F252').Value = 0.7392142031453177
$ws.Range('G252').Value = 1.379606495951268
$ws.Range('H252').Value = 'TZP'
$ws.Range('A253').Value = 'Age_cat0–28 d'
$ws.Range('B253').Value = 0.7452352981706942
$ws.Range('C253').Value = 0.3539073716208163
$ws.Range('D253').Value = -0.8308820261385873
$ws.Range('E253').Value = 0.4060402783715337
$ws.Range('F253').Value = 0.3690120214198593
$ws.Range('G253').Value = 1.481826769133568
$ws.Range('H253').Value = 'TZP'
$ws.Range('A254').Value = 'Age_cat29–365 d'
$ws.Range('B254').Value = 0.8251223013987901
$ws.Range('C254').Value = 0.3368775103312124
$ws.Range('D254').Value = -0.5706040136837662
$ws.Range('E254').Value = 0.5682680971805545
$ws.Range('F254').Value = 0.4232036081975295
$ws.Range('G254').Value = 1.589184735589086
$ws.Range('H254').Value = 'TZP'
$ws.Range('A255').Value = 'Age_cat1–5 y'
$ws.Range('B255').Value = 1.022969034968856
$ws.Range('C255').Value = 0.3560010526225862
$ws.Range('D255').Value = 0.06378974864382075
$ws.Range('E255').Value = 0.9491376410321533
$ws.Range('F255').Value = 0.5043603389987948
$ws.Range('G255').Value = 2.043447546845866
$ws.Range('H255').Value = 'TZP'
$ws.Range('A256').Value = 'Age_cat6–<30 y'
$ws.Range('B256').Value = 1.241191349676181
$ws.Range('C256').Value = 0.2792237254772341
$ws.Range('D256').Value = 0.7738299597407585
$ws.Range('E256').Value = 0.4390313529777314
$ws.Range('F256').Value = 0.7172569937819174
$ws.Range('G256').Value = 2.146928216985859
$ws.Range('H256').Value = 'TZP'
$ws.Range('A257').Value = 'Age_cat52–<67 y'
$ws.Range('B257').Value = 0.6915584825649634
$ws.Range('C257').Value = 0.2644920459833628
$ws.Range('D257').Value = -1.394399429280131
$ws.Range('E257').Value = 0.1631970203610393
$ws.Range('F257').Value = 0.4101922704414107
$ws.Range('G257').Value = 1.159307971801076
$ws.Range('H257').Value = 'TZP'
$ws.Range('A258').Value = 'Age_cat≥67 y'
$ws.Range('B258').Value = 0.8674602363730076
$ws.Range('C258').Value = 0.2646451827060015
$ws.Range('D258').Value = -0.5372688208106362
$ws.Range('E258').Value = 0.5910819392114427
$ws.Range('F258').Value = 0.5153305044320108
$ws.Range('G258').Value = 1.457066074456465
$ws.Range('H258').Value = 'TZP'
$ws.Range('A259').Value = 'HospitalCHBH:Ward_ED_ICUED'
$ws.Range('B259').Value = 2.773613282477063
$ws.Range('C259').Value = 0.7842069013599707
$ws.Range('D259').Value = 1.300869582379733
$ws.Range('E259').Value = 0.1933030993888284
$ws.Range('F259').Value = 0.5130549673319463
$ws.Range('G259').Value = 11.92216803660624
$ws.Range('H259').Value = 'TZP'
$ws.Range('A260').Value = 'HospitalCNGMO:Ward_ED_ICUED'
$ws.Range('H260').Value = 'TZP'
$ws.Range('A261').Value = 'HospitalCHBH:Ward_ED_ICUOther'
$ws.Range('B261').Value = 3.919699389786237
$ws.Range('C261').Value = 0.4623546430404585
$ws.Range('D261').Value = 2.954474417299782
$ws.Range('E261').Value = 0.003132021197257005
$ws.Range('F261').Value = 1.603166091905406
$ws.Range('G261').Value = 9.863250856107802
$ws.Range('H261').Value = 'TZP'
$ws.Range('A262').Value = 'HospitalCNGMO:Ward_ED_ICUOther'
$ws.Range('H262').Value = 'TZP'
